$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "(according to the population census data)" row entirely - the
# exported table no longer carries a census-data qualifier.
$ws.Rows.Item(2).Delete()

# Keep only the 2014 figures - remove the old 1989/2002 columns (B:C); the
# surviving 2014 column slides left into column B.
$ws.Range("B1:C1").EntireColumn.Delete()

# The year-header cell no longer sits at the right edge of the table, so its
# right edge goes from a thick (medium) rule to a thin one, matching its new,
# interior position.
$ws.Range("B4").Borders.Item(10).LineStyle = 1
$ws.Range("B4").Borders.Item(10).Weight = 2

# Every row in the rebuilt table uses the taller 20.1pt row height.
$ws.Range("A1:A7").RowHeight = 20.1

Write-Output "done"
